# Apply the edit described by the diff to testdatalogin.xlsx
#
# Summary of the change on sheet "testdata":
#  - The old "Status / status / status" + "Pass / pass / text N" block that lived
#    in columns G:I is cleared out (content removed, row-1 formatting kept).
#  - A new "status" / "pass" column is written starting at column J.
#  - Row 1's "status" value is then filled right across J1:S1 (10 cells).
#  - Rows 2-5 only get a single "pass" value in column J.
#  - The unused "Status"/"Pass"/"text 2".."text 5" strings disappear from the
#    shared string table automatically because nothing references them anymore.
#  - Final selection ends up on C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata")

# Clear the old status/pass columns (keeps row 1 cell formatting, like a Delete keypress)
$ws.Range("G1:I5").ClearContents()

# Write the new status/pass column starting at J
$ws.Range("J1").Value = "status"
$ws.Range("J2").Value = "pass"
$ws.Range("J3").Value = "pass"
$ws.Range("J4").Value = "pass"
$ws.Range("J5").Value = "pass"

# Fill the header value right across the new columns J1:S1
$ws.Range("J1:S1").FillRight()

# Leave the selection where the author's cursor ended up
$ws.Range("C6").Select()
